$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Sequence"

$sequences = @(1,1,1,1,1,1,2,2,2,2,2,2,3,3,3,3,4,4,4,4)
for ($i = 0; $i -lt $sequences.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $sequences[$i]
}

$ws.Columns.Item(4).ColumnWidth = 11.6

$ws.Range("D6").Select()
